$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '28.868.56'
$ws.Range('D2').Style = "Normal"

$ws.Range('E3').Value = '  -0.60%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.819.50'
$ws.Range('D3').Style = "Normal"

$ws.Range('E4').Value = '  -0.75%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9925'
$ws.Range('D4').Style = "Normal"

$ws.Range('E5').Value = '  +0.44%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '242.49'
$ws.Range('D5').Style = "Normal"

$ws.Range('E6').Value = '  +0.52%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.6280'
$ws.Range('D6').Style = "Normal"

$ws.Range('E7').Value = '  -0.53%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.9955'
$ws.Range('D7').Style = "Normal"

$ws.Range('E8').Value = '  -1.75%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.07446'
$ws.Range('D8').Style = "Normal"

$ws.Range('E9').Value = '  +0.35%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.2920'
$ws.Range('D9').Style = "Normal"

$ws.Range('E10').Value = '  +1.20%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '22.94'
$ws.Range('D10').Style = "Normal"

$ws.Range('E11').Value = '  -1.22%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07660'
$ws.Range('D11').Style = "Normal"

$ws.Range('E12').Value = '  -0.76%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.814.71'
$ws.Range('D12').Style = "Normal"

$ws.Range('E13').Value = '  +0.50%  '

$ws.Range('E14').Value = '  +0.55%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.6646'
$ws.Range('D14').Style = "Normal"

$ws.Range('E15').Value = '  +0.54%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '82.69'
$ws.Range('D15').Style = "Normal"

$ws.Range('E16').Value = '  +2.67%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.000009660'
$ws.Range('D16').Style = "Normal"

$ws.Range('E17').Value = '  +0.73%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '6.004'
$ws.Range('D17').Style = "Normal"

$ws.Range('E18').Value = '  -0.36%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '28.879.74'
$ws.Range('D18').Style = "Normal"

$ws.Range('E19').Value = '  +1.74%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.52'
$ws.Range('D19').Style = "Normal"

$ws.Range('E20').Value = '  -0.33%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '223.09'
$ws.Range('D20').Style = "Normal"

$ws.Range('E21').Value = '  -0.60%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.9948'
$ws.Range('D21').Style = "Normal"

$ws.Range('E22').Value = '  -1.41%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.089'
$ws.Range('D22').Style = "Normal"

$ws.Range('E23').Value = '  -0.60%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.9955'
$ws.Range('D23').Style = "Normal"

$ws.Range('E24').Value = '  -0.19%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '158.97'
$ws.Range('D24').Style = "Normal"

$ws.Range('E25').Value = '  +3.22%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.1404'
$ws.Range('D25').Style = "Normal"

$ws.Range('E26').Value = '  +0.55%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.450'
$ws.Range('D26').Style = "Normal"

$ws.Range('E27').Value = '  +0.02%  '

$ws.Range('E28').Value = '  -0.11%  '

$ws.Range('E29').Value = '  +1.18%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '4.100'
$ws.Range('D29').Style = "Normal"

$ws.Range('E30').Value = '  +0.46%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.035'
$ws.Range('D30').Style = "Normal"

$ws.Range('E31').Value = '  +4.92%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.05437'
$ws.Range('D31').Style = "Normal"

$ws.Range('E32').Value = '  -0.31%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.194'
$ws.Range('D32').Style = "Normal"

$ws.Range('E33').Value = '  +0.20%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.847'
$ws.Range('D33').Style = "Normal"

$ws.Range('E34').Value = '  +0.45%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7388'
$ws.Range('D34').Style = "Normal"

$ws.Range('E35').Value = '  -1.15%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.130'
$ws.Range('D35').Style = "Normal"

$ws.Range('E36').Value = '  -3.51%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.605'
$ws.Range('D36').Style = "Normal"

$ws.Range('E37').Value = '  -2.28%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.229.63'
$ws.Range('D37').Style = "Normal"

$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('E38').Value = '  -0.25%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01778'
$ws.Range('D38').Style = "Normal"

$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('E39').Value = '  -0.81%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.738'
$ws.Range('D39').Style = "Normal"

$ws.Range('E40').Value = '  +6.16%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '6.642'
$ws.Range('D40').Style = "Normal"

$ws.Range('E41').Value = '  +0.80%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.8932'
$ws.Range('D41').Style = "Normal"

$ws.Range('E42').Value = '  -0.58%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.9952'
$ws.Range('D42').Style = "Normal"

$ws.Range('E43').Value = '  -0.31%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '101.09'
$ws.Range('D43').Style = "Normal"

$ws.Range('B44').Value = 'BabyDogeCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.00000000123'
$ws.Range('D44').Style = "Normal"

$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('E45').Value = '  +1.11%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '64.94'
$ws.Range('D45').Style = "Normal"

$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('E46').Value = '  -1.13%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.5054'
$ws.Range('D46').Style = "Normal"

$ws.Range('B47').Value = 'TheSandbox'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('E47').Value = '  +1.41%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.4026'
$ws.Range('D47').Style = "Normal"

$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E48').Value = '  +1.53%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '8.972'
$ws.Range('D48').Style = "Normal"

$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('E49').Value = '  +1.93%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.657'
$ws.Range('D49').Style = "Normal"

$ws.Range('B50').Value = 'XinFinNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range('E50').Value = '  +2.02%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.07191'
$ws.Range('D50').Style = "Normal"

$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('E51').Value = '  +0.60%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05778'
$ws.Range('D51').Style = "Normal"
